# Apply the "Add files via upload" edit:
#  - new redirect-check table (K1:N15)
#  - new snippet notes (F17:G19) and a "redirect exists" note (A52)
#  - a new competitor-analysis table (A55:E57)
#  - assorted header/column formatting to support the above

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Column widths (new/changed columns)
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth  = 20.5   # C: was 13.14 -> ~21.28
$ws.Columns.Item(4).ColumnWidth  = 19.17  # D: new -> 20
$ws.Columns.Item(5).ColumnWidth  = 11.33  # E: new -> 12.14
$ws.Columns.Item(6).ColumnWidth  = 12.33  # F: new -> 13.14
$ws.Columns.Item(11).ColumnWidth = 10     # K: new -> 10.86
$ws.Columns.Item(12).ColumnWidth = 37.33  # L: new -> 38.14
$ws.Columns.Item(13).ColumnWidth = 10     # M: new -> 10.86
$ws.Columns.Item(14).ColumnWidth = 12.17  # N: new -> 13

# ---------------------------------------------------------------------------
# 2. Header row 1 - bold/size-14 styling for the existing headers, plus the
#    four new headers for the redirect-check table
# ---------------------------------------------------------------------------
$headerRange = $ws.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.Font.Size = 14

$ws.Range("K1").Value = "#"
$ws.Range("L1").Value = "url"
$ws.Range("M1").Value = "Ответ"
$ws.Range("N1").Value = "Результат"
$newHeaderRange = $ws.Range("K1:N1")
$newHeaderRange.Font.Bold = $true
$newHeaderRange.Font.Size = 14

# ---------------------------------------------------------------------------
# 3. Redirect-check table, rows 2-15
# ---------------------------------------------------------------------------
$urls = @(
  "https://nowblepro.github.io/index.html",
  "https://nowblepro.github.io/",
  "https://nowblepro.github.io//",
  "https://nowblepro.github.io///",
  "https://nowblepro.github.io////",
  "https://nowblepro.github.io/////",
  "https://nowblepro.github.io//////",
  "https://nowblepro.github.io///////",
  "https://nowblepro.github.io////////",
  "https://nowblepro.github.io/////////",
  "https://nowblepro.github.io//////////",
  "https://nowblepro.github.io/?",
  "https://nowblepro.github.io/index",
  "https://nowblepro.github.io"
)

for ($i = 0; $i -lt $urls.Length; $i++) {
  $row = 2 + $i
  $ws.Cells.Item($row, 11).Value = $i + 1
  $ws.Cells.Item($row, 12).Value = $urls[$i]
  $ws.Cells.Item($row, 13).Value = 200
  $ws.Cells.Item($row, 14).Value = "доступна"
}

# ---------------------------------------------------------------------------
# 4. Snippet notes, rows 17-19 (columns F/G)
# ---------------------------------------------------------------------------
$ws.Range("F17").Value = "Сниппеты"
$ws.Range("F18").Value = "title"
$ws.Range("F19").Value = "description"
$snippetLabels = $ws.Range("F17:F19")
$snippetLabels.Font.Bold = $true
$snippetLabels.Font.Size = 12

$ws.Range("G18").Value = "<title> Купить одежду для женщин, мужчин, детей и военных в спб 😀 [лучише цены] </title>"
$ws.Range("G19").Value = "<description>👦👧 Одежда на любой вкус в спб для детей и их родителей. 😎 Так же у нас есть верхняя одежда. У нас намного дешевле. Одежда для военных или же тактическая. 😎 Женская и мужская одежда на любое время года .</description>"

# ---------------------------------------------------------------------------
# 5. Redirect note, row 52
# ---------------------------------------------------------------------------
$ws.Range("A52").Value = "Есть переадресация с //, /, /index.html"

# ---------------------------------------------------------------------------
# 6. Competitor-analysis table, rows 55-57
# ---------------------------------------------------------------------------
$ws.Range("A55").Value = "Конкуренты"
$ws.Range("B55").Value = "Минусы"
$ws.Range("C55").Value = "Решение"
$ws.Range("D55").Value = "Интерактивность"
$ws.Range("E55").Value = "Вывод"

$ws.Range("A55").Font.Bold = $true
$ws.Range("A55").Font.Size = 16
$ws.Range("B55:E55").Font.Bold = $true
$ws.Range("B55:E55").Font.Size = 11

$ws.Range("A56").Value = "hitchhikers.ru"
$ws.Range("B56").Value = "Не равномерный размер блоков, из-за этого сайт не подстраивается под разные разрешения мониторов, поэтому все выглядит неуклюже."
$ws.Range("C56").Value = "Сделать равномерные блоки. Адаптировать сайт под разные мониторы и их разрешения."
$ws.Range("D56").Value = "Интерактивности мало, сайт минималистичен, присутствует слайдер небольшой"
$ws.Range("E56").Value = "Интерактивновти мало, нужно дороботать"

$ws.Range("A57").Value = "https://mynamestore.ru/"
$ws.Range("B57").Value = "Не подходящий цвет шрифта относительно фона. Видео плеер не должен появляться. При нажатие на картинки выводит на не существующую страницу. Нет страницы ERROR 404."
$ws.Range("C57").Value = "Определить самый благоприятный для глаза цвет шрифта относительно фона. Сделать плеер не интерактивным для гостя. "
$ws.Range("D57").Value = "Интерактивности нет, сайт практически одностраничник"
$ws.Range("E57").Value = "Добавить интерактив, слайдер как минимум"

# Alignment for rows 56-57
$ws.Range("A56:A57").HorizontalAlignment = -4108  # xlCenter
$ws.Range("A56:A57").VerticalAlignment   = -4108  # xlCenter

$ws.Range("B56").WrapText = $true
$ws.Range("B56").VerticalAlignment = -4160        # xlTop
$ws.Range("D56:E56").WrapText = $true
$ws.Range("D56:E56").VerticalAlignment = -4160
$ws.Range("D57:E57").WrapText = $true
$ws.Range("D57:E57").VerticalAlignment = -4160

$ws.Range("C56").WrapText = $true
$ws.Range("C56").HorizontalAlignment = -4131      # xlLeft
$ws.Range("C56").VerticalAlignment = -4160

$ws.Range("B57:C57").WrapText = $true
$ws.Range("B57:C57").HorizontalAlignment = -4131
$ws.Range("B57:C57").VerticalAlignment = -4160

$ws.Rows.Item(56).RowHeight = 144.75
$ws.Rows.Item(57).RowHeight = 158.25

# ---------------------------------------------------------------------------
# 7. View settings: zoom + selection
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.Zoom = 85
$ws.Range("F57").Select()
